$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PlanningProjet")
$ws.Activate()

$win = $excel.ActiveWindow
$win.Zoom = 54

$ws.Range("C25").Select()

$ws.Columns.Item(2).ColumnWidth = 19.04

$ws.Rows.Item(25).RowHeight = 30
